$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.226.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.446.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.18%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.107"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.361"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.26%  "
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.884.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.095.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.82%  "
$ws.Range("E16").Value = "  +5.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.455.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.05%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0798"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  +6.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.418"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "316.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "144.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0965"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0527"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("B47").Value = "Polygon"
$ws.Range("C47").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.411"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.84%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.575"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("E49").Value = "  +3.05%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  +4.94%  "
